$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff inserts one brand-new data row at sheet row 42 ("1a (guarda)"
# Asterix reading dated 44469), which pushes every existing row from
# 42..115 down by one (43..116). Using Rows.Insert reproduces that shift
# for every column (including cell styles, e.g. the date-formatted D
# column) without having to rewrite the rest of the sheet by hand.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with its data.
$ws.Cells.Item(42, 1).Value = 11
$ws.Cells.Item(42, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value = "Bíobío"
$ws.Cells.Item(42, 4).Value = 44469
$ws.Cells.Item(42, 5).Value = 8
$ws.Cells.Item(42, 6).Value = 100114001
$ws.Cells.Item(42, 7).Value = "Papa"
$ws.Cells.Item(42, 8).Value = "Asterix"
$ws.Cells.Item(42, 9).Value = "1a (guarda)"
$ws.Cells.Item(42, 10).Value = 2000
$ws.Cells.Item(42, 11).Value = 9500
$ws.Cells.Item(42, 12).Value = 10000
$ws.Cells.Item(42, 13).Value = 9750
$ws.Cells.Item(42, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Arauco"
$ws.Cells.Item(42, 16).Value = 390
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
